# Applies the "Add files via upload" edit: appends 43 new BPR/BPRS entries
# (columns B "BankId" and C "BankName" only, per the source diff) to rows
# 188-230 of Sheet1, and fixes up the _FilterDatabase defined name + the
# worksheet AutoFilter range (both stay anchored at row 187, the end of the
# original filtered table) plus the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-apply the AutoFilter over the ORIGINAL table extent (A1:D187) first,
# while no data exists below row 187, so Excel does not silently grow the
# filter range to cover the rows we are about to add.
$ws.Range("A1:D1").AutoFilter()
$ws.Range("A1:D187").AutoFilter()

# --- Fix the hidden _xlnm._FilterDatabase defined name to match.
$filterName = $wb.Names.Item("Sheet1!_FilterDatabase")
$filterName.RefersTo = "=Sheet1!`$A`$1:`$D`$187"

# --- Append the new BankId / BankName rows (188-230). Columns A (BankId)
# and D (effective date) are intentionally left blank, matching the source.
$ws.Range("B188").Value = 600073
$ws.Range("C188").Value = "PT BPR Mega Karsa Mandiri"
$ws.Range("B189").Value = 600100
$ws.Range("C189").Value = "PT BPR Cinere Artha Raya"
$ws.Range("B190").Value = 600101
$ws.Range("C190").Value = "PT BPR Sisibahari Dana"
$ws.Range("B191").Value = 600150
$ws.Range("C191").Value = "PT BPR Bintang Ekonomi Sejahtera"
$ws.Range("B192").Value = 600159
$ws.Range("C192").Value = "PT BPR Nusa Galang Makmur"
$ws.Range("B193").Value = 600266
$ws.Range("C193").Value = "PT BPR EDCCASH"
$ws.Range("B194").Value = 600775
$ws.Range("C194").Value = "PT BPR Kudamas Sentosa"
$ws.Range("B195").Value = 600797
$ws.Range("C195").Value = "PT BPR Mustika Utama Kolaka"
$ws.Range("B196").Value = 600850
$ws.Range("C196").Value = "PT BPR Dhasatra Artha Sempurna"
$ws.Range("B197").Value = 600861
$ws.Range("C197").Value = "PT BPR Nova Trijaya"
$ws.Range("B198").Value = 600920
$ws.Range("C198").Value = "PT BPR Pancadana"
$ws.Range("B199").Value = 601014
$ws.Range("C199").Value = "PT BPR Legian"
$ws.Range("B200").Value = 601093
$ws.Range("C200").Value = "PT BPR Aceh Utara"
$ws.Range("B201").Value = 601156
$ws.Range("C201").Value = "PT BPR LPN Kampung Baru"
$ws.Range("B202").Value = 601160
$ws.Range("C202").Value = "PT BPR Cahaya Nagari"
$ws.Range("B203").Value = 601182
$ws.Range("C203").Value = "PT BPR Carano Nagari"
$ws.Range("B204").Value = 601212
$ws.Range("C204").Value = "BPR LPN Kampung Manggis"
$ws.Range("B205").Value = 601227
$ws.Range("C205").Value = "PT BPR Mitra Danagung"
$ws.Range("B206").Value = 601282
$ws.Range("C206").Value = "PT BPR Arthasraya Sejahtera"
$ws.Range("B207").Value = 601298
$ws.Range("C207").Value = "PT BPR Bina Dian Citra"
$ws.Range("B208").Value = 601304
$ws.Range("C208").Value = "PT BPR Cita Makmur Lestari"
$ws.Range("B209").Value = 601318
$ws.Range("C209").Value = "PT BPR Lumasindo Perkasa Putra"
$ws.Range("B210").Value = 601327
$ws.Range("C210").Value = "PT BPR Kop. Jawa Barat"
$ws.Range("B211").Value = 601344
$ws.Range("C211").Value = "PT BPR Mutiara Artha Pratama"
$ws.Range("B212").Value = 601473
$ws.Range("C212").Value = "PT BPR Sinar Baru Perkasa"
$ws.Range("B213").Value = 601741
$ws.Range("C213").Value = "PT BPR Artha Dharma"
$ws.Range("B214").Value = 601790
$ws.Range("C214").Value = "PT BPR Cakra Dharma Artamandiri"
$ws.Range("B215").Value = 601802
$ws.Range("C215").Value = "PT BPR Kujang Artha Sembada"
$ws.Range("B216").Value = 601836
$ws.Range("C216").Value = "PT BPR KS Bali Agung Sedana"
$ws.Range("B217").Value = 601880
$ws.Range("C217").Value = "PT BPR Tugu Kencana"
$ws.Range("B218").Value = 601920
$ws.Range("C218").Value = "PT BPR Akarumi"
$ws.Range("B219").Value = 601980
$ws.Range("C219").Value = "PT BPR Agra Arthaka Mulya"
$ws.Range("B220").Value = 601989
$ws.Range("C220").Value = "PT BPR Sambas Arta"
$ws.Range("B221").Value = 602013
$ws.Range("C221").Value = "PT BPR Mitra Bunda Mandiri"
$ws.Range("B222").Value = 602060
$ws.Range("C222").Value = "PT BPR Sinarenam Permai Jatiasih"
$ws.Range("B223").Value = 602545
$ws.Range("C223").Value = "PT BPR Budisetia"
$ws.Range("B224").Value = 602589
$ws.Range("C224").Value = "PT BPR Indomitra Mega Kapital"
$ws.Range("B225").Value = 602603
$ws.Range("C225").Value = "PT BPR Vox Modern Danamitra"
$ws.Range("B226").Value = 602634
$ws.Range("C226").Value = "PT BPR Bungo Mandiri"
$ws.Range("B227").Value = 620072
$ws.Range("C227").Value = "PT BPRS Al-Hidayah"
$ws.Range("B228").Value = 620090
$ws.Range("C228").Value = " PT BPRS Safir Bengkulu"
$ws.Range("B229").Value = 620103
$ws.Range("C229").Value = "PT BPRS Jabal Tsur"
$ws.Range("B230").Value = 620119
$ws.Range("C230").Value = "PT BPRS Hidayah"

# --- Match the saved selection/scroll state from the edited workbook.
$ws.Range("A188").Select()

